# planilha-dinamica.xlsx - "adicionando projeto sobre tratamento, analise e
# insights dos dados com inteligencia artificial"
#
# The commit unhides the three helper sheets (Data, Economias, Controller)
# that back the dashboard's pivot tables/slicers, and switches the
# selected/active sheet from "Dashboard" to "Controller" so the workbook
# opens showing the controller sheet instead.

$wb = $excel.ActiveWorkbook

# Unhide the previously hidden worksheets.
$wb.Worksheets.Item("Data").Visible = -1
$wb.Worksheets.Item("Economias").Visible = -1
$wb.Worksheets.Item("Controller").Visible = -1

# Make "Controller" the active/selected sheet (was "Dashboard").
$wb.Worksheets.Item("Controller").Activate()
